$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The match "Always Ready" vs "Bolivar" (row with Id Q94GGB6s, Time 17:00)
# was removed from the weekly fixtures sheet. Delete its entire row so the
# rows below it shift up.
$ws.Rows.Item(2).Delete()
